$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-25 Saturday" "2025-10-26 Sunday"

Replace-Text "615×6=" "154×6="
Replace-Text "597×9=" "483×9="
Replace-Text "469×2=" "537×7="
Replace-Text "256×3=" "660×6="
Replace-Text "128×2=" "914×4="

Replace-Text "592×4=" "843×2="
Replace-Text "790×4=" "661×6="
Replace-Text "951×4=" "837×4="
Replace-Text "388×9=" "230×7="
Replace-Text "512×7=" "123×9="

Replace-Text "783×3=" "856×6="
Replace-Text "195×3=" "137×3="
Replace-Text "760×9=" "458×9="
Replace-Text "156×8=" "615×6="
Replace-Text "332×8=" "510×4="

Replace-Text "700×3=" "660×5="
Replace-Text "479×2=" "642×9="
Replace-Text "219×7=" "857×5="
Replace-Text "720×5=" "117×8="
Replace-Text "873×8=" "447×2="

Replace-Text "619×3=" "665×2="
Replace-Text "236×3=" "582×8="
Replace-Text "718×6=" "910×6="
Replace-Text "648×5=" "416×4="
Replace-Text "649×8=" "663×3="
